$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text (volume number, report week date range) ---
$ws.Range("A8").Value = "Volume 30   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/30/2023  Through  11/5/2023"

# --- Donor cells used to copy the "placeholder" text formatting (style 14)
#     onto cells that transition from a numeric value to a text placeholder ---
$donorZero = $ws.Cells.Item(14,4)   # D14: text "0", style 14
$donorStar = $ws.Cells.Item(22,14)  # N22: text "***.*", style 14

$ws.Cells.Item(15,4).NumberFormat = '#,##0'
$ws.Cells.Item(15,4).Value = 2
$ws.Cells.Item(15,5).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(15,5).Value = -100
$ws.Cells.Item(15,7).Value = 3
$ws.Cells.Item(15,8).Value = -33.333333333333
$ws.Cells.Item(15,10).Value = 13
$ws.Cells.Item(15,11).Value = -30.769230769230
$ws.Cells.Item(15,12).Value = 0
$ws.Cells.Item(15,13).Value = -43.75
$ws.Cells.Item(15,14).Value = -72.727272727272
$ws.Cells.Item(16,3).Value = 3
$ws.Cells.Item(16,4).Value = 3
$ws.Cells.Item(16,5).Value = 0
$ws.Cells.Item(16,6).Value = 15
$ws.Cells.Item(16,7).Value = 14
$ws.Cells.Item(16,8).Value = 7.142857142857
$ws.Cells.Item(16,9).Value = 145
$ws.Cells.Item(16,10).Value = 158
$ws.Cells.Item(16,11).Value = -8.227848101265
$ws.Cells.Item(16,12).Value = 68.604651162790
$ws.Cells.Item(16,13).Value = -38.297872340425
$ws.Cells.Item(16,14).Value = -84.239130434782
$ws.Cells.Item(17,3).Value = 5
$ws.Cells.Item(17,4).Value = 8
$ws.Cells.Item(17,5).Value = -37.5
$ws.Cells.Item(17,6).Value = 26
$ws.Cells.Item(17,7).Value = 32
$ws.Cells.Item(17,8).Value = -18.75
$ws.Cells.Item(17,9).Value = 319
$ws.Cells.Item(17,10).Value = 328
$ws.Cells.Item(17,11).Value = -2.743902439024
$ws.Cells.Item(17,12).Value = 5.280528052805
$ws.Cells.Item(17,13).Value = 74.316939890710
$ws.Cells.Item(17,14).Value = -50.619195046439
$ws.Cells.Item(18,3).Value = 4
$ws.Cells.Item(18,4).Value = 4
$ws.Cells.Item(18,5).Value = 0
$ws.Cells.Item(18,6).Value = 9
$ws.Cells.Item(18,7).Value = 20
$ws.Cells.Item(18,8).Value = -55
$ws.Cells.Item(18,9).Value = 94
$ws.Cells.Item(18,10).Value = 143
$ws.Cells.Item(18,11).Value = -34.265734265734
$ws.Cells.Item(18,12).Value = -35.616438356164
$ws.Cells.Item(18,13).Value = -32.857142857142
$ws.Cells.Item(18,14).Value = -88.941176470588
$ws.Cells.Item(19,3).Value = 13
$ws.Cells.Item(19,4).Value = 14
$ws.Cells.Item(19,5).Value = -7.142857142857
$ws.Cells.Item(19,6).Value = 28
$ws.Cells.Item(19,7).Value = 38
$ws.Cells.Item(19,8).Value = -26.315789473684
$ws.Cells.Item(19,9).Value = 394
$ws.Cells.Item(19,10).Value = 454
$ws.Cells.Item(19,11).Value = -13.215859030837
$ws.Cells.Item(19,12).Value = 15.542521994134
$ws.Cells.Item(19,13).Value = -20.081135902636
$ws.Cells.Item(19,14).Value = -33.106960950764
$ws.Cells.Item(20,3).Value = 5
$ws.Cells.Item(20,4).NumberFormat = '#,##0'
$ws.Cells.Item(20,4).Value = 1
$ws.Cells.Item(20,5).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(20,5).Value = 400
$ws.Cells.Item(20,6).Value = 13
$ws.Cells.Item(20,7).Value = 3
$ws.Cells.Item(20,8).Value = 333.333333333333
$ws.Cells.Item(20,9).Value = 98
$ws.Cells.Item(20,10).Value = 104
$ws.Cells.Item(20,11).Value = -5.769230769230
$ws.Cells.Item(20,12).Value = 53.125
$ws.Cells.Item(20,13).Value = 1.030927835051
$ws.Cells.Item(20,14).Value = -90.268123138033
$ws.Cells.Item(21,3).Value = 30
$ws.Cells.Item(21,4).Value = 32
$ws.Cells.Item(21,5).Value = -6.25
$ws.Cells.Item(21,6).Value = 93
$ws.Cells.Item(21,7).Value = 110
$ws.Cells.Item(21,8).Value = -15.454545454545
$ws.Cells.Item(21,9).Value = 1064
$ws.Cells.Item(21,10).Value = 1205
$ws.Cells.Item(21,11).Value = -11.701244813278
$ws.Cells.Item(21,12).Value = 11.180773249738
$ws.Cells.Item(21,13).Value = -9.059829059829
$ws.Cells.Item(21,14).Value = -73.818897637795
$ws.Cells.Item(22,3).NumberFormat = '@'
$ws.Cells.Item(22,3).Value = "0"
$donorZero.Copy()
$ws.Cells.Item(22,3).PasteSpecial(-4122)
$ws.Cells.Item(22,4).NumberFormat = '#,##0'
$ws.Cells.Item(22,4).Value = 1
$ws.Cells.Item(22,5).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(22,5).Value = -100
$ws.Cells.Item(22,6).Value = 2
$ws.Cells.Item(22,7).Value = 3
$ws.Cells.Item(22,8).Value = -33.333333333333
$ws.Cells.Item(22,10).Value = 26
$ws.Cells.Item(22,11).Value = -38.461538461538
$ws.Cells.Item(23,3).Value = 1
$ws.Cells.Item(23,4).Value = 5
$ws.Cells.Item(23,5).Value = -80
$ws.Cells.Item(23,7).Value = 20
$ws.Cells.Item(23,8).Value = -35
$ws.Cells.Item(23,9).Value = 145
$ws.Cells.Item(23,10).Value = 151
$ws.Cells.Item(23,11).Value = -3.973509933774
$ws.Cells.Item(23,12).Value = -8.805031446540
$ws.Cells.Item(23,13).Value = 40.776699029126
$ws.Cells.Item(24,3).Value = 24
$ws.Cells.Item(24,4).Value = 25
$ws.Cells.Item(24,5).Value = -4
$ws.Cells.Item(24,6).Value = 72
$ws.Cells.Item(24,7).Value = 89
$ws.Cells.Item(24,8).Value = -19.101123595505
$ws.Cells.Item(24,9).Value = 855
$ws.Cells.Item(24,10).Value = 1053
$ws.Cells.Item(24,11).Value = -18.803418803418
$ws.Cells.Item(24,12).Value = 14.457831325301
$ws.Cells.Item(24,13).Value = -17.709335899903
$ws.Cells.Item(25,3).Value = 10
$ws.Cells.Item(25,4).Value = 6
$ws.Cells.Item(25,5).Value = 66.666666666666
$ws.Cells.Item(25,6).Value = 41
$ws.Cells.Item(25,7).Value = 39
$ws.Cells.Item(25,8).Value = 5.128205128205
$ws.Cells.Item(25,9).Value = 530
$ws.Cells.Item(25,10).Value = 485
$ws.Cells.Item(25,11).Value = 9.278350515463
$ws.Cells.Item(25,12).Value = 39.473684210526
$ws.Cells.Item(25,13).Value = 14.718614718614
$ws.Cells.Item(26,4).Value = 2
$ws.Cells.Item(26,7).Value = 5
$ws.Cells.Item(26,8).Value = -60
$ws.Cells.Item(26,10).Value = 25
$ws.Cells.Item(26,11).Value = -40
$ws.Cells.Item(26,12).Value = -25
$ws.Cells.Item(27,3).Value = 4
$ws.Cells.Item(27,4).NumberFormat = '@'
$ws.Cells.Item(27,4).Value = "0"
$donorZero.Copy()
$ws.Cells.Item(27,4).PasteSpecial(-4122)
$ws.Cells.Item(27,5).NumberFormat = '@'
$ws.Cells.Item(27,5).Value = "***.*"
$donorStar.Copy()
$ws.Cells.Item(27,5).PasteSpecial(-4122)
$ws.Cells.Item(27,6).Value = 7
$ws.Cells.Item(27,7).Value = 12
$ws.Cells.Item(27,8).Value = -41.666666666666
$ws.Cells.Item(27,9).Value = 48
$ws.Cells.Item(27,11).Value = -5.882352941176
$ws.Cells.Item(27,12).Value = 54.838709677419
$ws.Cells.Item(30,4).NumberFormat = '@'
$ws.Cells.Item(30,4).Value = "0"
$donorZero.Copy()
$ws.Cells.Item(30,4).PasteSpecial(-4122)
$ws.Cells.Item(30,5).NumberFormat = '@'
$ws.Cells.Item(30,5).Value = "***.*"
$donorStar.Copy()
$ws.Cells.Item(30,5).PasteSpecial(-4122)
